$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 2),
    @(1, 2),
    @(1, 3),
    @(1, 3),
    @(1, 3),
    @(1, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 13 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("B9").Select()
